$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.697.99"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "3.800.92"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  +0.06%  "
$c = $ws.Range("D5")
$c.Value = "'602.63"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.27%  "
$c = $ws.Range("D6")
$c.Value = "'165.72"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("E7").Value = "  +0.12%  "
$c = $ws.Range("D8")
$c.Value = "'0.518"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("E11").Value = "  -0.33%  "
$c = $ws.Range("D12")
$c.Value = "'0.0000250"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.45%  "
$c = $ws.Range("D13")
$c.Value = "'35.78"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").Value = "4.445.36"
$ws.Range("E14").Value = "  +0.88%  "
$ws.Range("D15").Value = "3.811.42"
$ws.Range("E15").Value = "  +1.17%  "
$c = $ws.Range("D16")
$c.Value = "'18.46"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.97%  "
$ws.Range("D17").Value = "67.724.12"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("E19").Value = "  +1.41%  "
$c = $ws.Range("D20")
$c.Value = "'462.95"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.40%  "
$c = $ws.Range("D21")
$c.Value = "'9.85"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.63%  "
$c = $ws.Range("D22")
$c.Value = "'0.699"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("E23").Value = "  -3.65%  "
$c = $ws.Range("D24")
$c.Value = "'83.17"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.01%  "
$c = $ws.Range("D25")
$c.Value = "'12.09"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.63%  "
$c = $ws.Range("D26")
$c.Value = "'2.11"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.77%  "
$c = $ws.Range("D27")
$c.Value = "'10.04"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").Value = "3.951.86"
$ws.Range("E29").Value = "  +0.78%  "
$c = $ws.Range("D30")
$c.Value = "'2.78"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.79%  "
$c = $ws.Range("D31")
$c.Value = "'7.34"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.36%  "
$c = $ws.Range("D32")
$c.Value = "'2.21"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.26%  "
$c = $ws.Range("D33")
$c.Value = "'29.44"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("E34").Value = "  +0.10%  "
$c = $ws.Range("D35")
$c.Value = "'9.07"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.84%  "
$c = $ws.Range("D36")
$c.Value = "'0.0995"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D38")
$c.Value = "'0.997"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D39")
$c.Value = "'5.80"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D40")
$c.Value = "'3.23"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.90%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("E42").Value = "  +0.05%  "
$c = $ws.Range("D43")
$c.Value = "'44.57"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -4.19%  "
$c = $ws.Range("D44")
$c.Value = "'47.64"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.28%  "
$c = $ws.Range("D45")
$c.Value = "'0.299"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.09%  "
$c = $ws.Range("D46")
$c.Value = "'28.39"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +9.28%  "
$c = $ws.Range("D47")
$c.Value = "'1.40"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +13.00%  "
$c = $ws.Range("D48")
$c.Value = "'150.76"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.16%  "
$c = $ws.Range("D49")
$c.Value = "'8.34"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("E50").Value = "  +1.50%  "
$c = $ws.Range("D51")
$c.Value = "'388.30"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.30%  "
